$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: was "salary" -> now "Dad", amount unchanged, date updated
$ws.Range("A2").Value = "Dad"
$ws.Range("B2").Value = 100000
$ws.Range("C2").Value = 45978.22928240741

# Row 3: new row "Job"
$ws.Range("A3").Value = "Job"
$ws.Range("B3").Value = 1000000
$ws.Range("C3").Value = 45962.22928240741

# Row 4: new row "interest"
$ws.Range("A4").Value = "interest"
$ws.Range("B4").Value = 100
$ws.Range("C4").Value = 45946.22928240741

# Row 5: original "salary" row, shifted down to row 5
$ws.Range("A5").Value = "salary"
$ws.Range("B5").Value = 100000
$ws.Range("C5").Value = 45658.22928240741

# Row 6: new row "rent"
$ws.Range("A6").Value = "rent"
$ws.Range("B6").Value = 100000
$ws.Range("C6").Value = 45658.22928240741

# Row 7: new row "real-estate"
$ws.Range("A7").Value = "real-estate"
$ws.Range("B7").Value = 12000
$ws.Range("C7").Value = 45658.22928240741

# Copy the date number format (numFmtId 14) from C2 onto the rest of the date
# column so all new date cells render/save using the same style.
$ws.Range("C2").Copy()
$ws.Range("C3:C7").PasteSpecial(-4122)
$excel.CutCopyMode = $false
